# Architectures of the internet of things.pptx
# - Add body text to slide 2 ("Internet enabled smart devices")
# - Append 4 new "Title and Content" slides (Methodology, Features,
#   Demonstration of product, Next steps)

$p = $ppt.ActivePresentation
$layout = $p.SlideMaster.CustomLayouts.Item(2)   # "Title and Content"

# ---------------------------------------------------------------------------
# Slide 2: fill in the previously-empty content placeholder.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$s2body = $s2.Shapes.Item(2).TextFrame.TextRange
$s2body.Text = "An Internet enabled smart device is a device that tends to be connected other smart devices and/or the internet via various wireless protocols such as WI-FI, LI-FI, NFC etc."
$s2body.LanguageID = "en-GB"
$s2body.InsertAfter("`rThe smart device in this presentation comprises of a Raspberry Pi 3 running Raspbian operating system connected to the IBM Bluemix platform for data capture and processing.")
$s2.Shapes.Item(2).TextFrame.TextRange.LanguageID = "en-GB"

# ---------------------------------------------------------------------------
# Slide 3: Methodology
# ---------------------------------------------------------------------------
$s3 = $p.Slides.AddSlide(3, $layout)

$s3title = $s3.Shapes.Item(1).TextFrame.TextRange
$s3title.Text = "Methodology"
$s3title.LanguageID = "en-GB"

$s3body = $s3.Shapes.Item(2).TextFrame.TextRange
$s3body.Text = "Setup Raspberry Pi for use"
$s3body.LanguageID = "en-GB"
$s3body.InsertAfter("`rInstalled needed software")
$s3body.InsertAfter("`rNode-red")
$s3body.InsertAfter("`rConnected Node-red to IBM Bluemix platform for data collection and processing")
$s3body.InsertAfter("`rSetup IBM Bluemix platform for data processing")
$s3bodyAll = $s3.Shapes.Item(2).TextFrame.TextRange
$s3bodyAll.LanguageID = "en-GB"
$s3bodyAll.Paragraphs(3,1).IndentLevel = 2

# ---------------------------------------------------------------------------
# Slide 4: Features
# ---------------------------------------------------------------------------
$s4 = $p.Slides.AddSlide(4, $layout)

$s4title = $s4.Shapes.Item(1).TextFrame.TextRange
$s4title.Text = "Features"
$s4title.LanguageID = "en-GB"

$s4body = $s4.Shapes.Item(2).TextFrame.TextRange
$s4body.Text = "CPU temp of RPI is sent to IBM Bluemix and processed into a real time updating graph."
$s4body.LanguageID = "en-GB"
$s4body.InsertAfter("`rCode has been implemented for an Arduino system to use a serial connection to get data to the RPI however this is untested as IBM Bluemix is unavailable.")
$s4.Shapes.Item(2).TextFrame.TextRange.LanguageID = "en-GB"

# ---------------------------------------------------------------------------
# Slide 5: Demonstration of product
# ---------------------------------------------------------------------------
$s5 = $p.Slides.AddSlide(5, $layout)

$s5title = $s5.Shapes.Item(1).TextFrame.TextRange
$s5title.Text = "Demonstration of product"
$s5title.LanguageID = "en-GB"

# ---------------------------------------------------------------------------
# Slide 6: Next steps
# ---------------------------------------------------------------------------
$s6 = $p.Slides.AddSlide(6, $layout)

$s6title = $s6.Shapes.Item(1).TextFrame.TextRange
$s6title.Text = "Next steps"
$s6title.LanguageID = "en-GB"

$s6body = $s6.Shapes.Item(2).TextFrame.TextRange
$s6body.Text = "Connecting the system to other internet enabled devices that use either WI-FI or LI-FI so that other areas of buildings can be monitored."
$s6body.LanguageID = "en-GB"
$s6body.InsertAfter("`rDevices would be fitted to monitor temperature humidity and light level.")
$s6bodyAll = $s6.Shapes.Item(2).TextFrame.TextRange
$s6bodyAll.LanguageID = "en-GB"
$s6bodyAll.Paragraphs(2,1).IndentLevel = 2
